# Reworks the "add" unit-test sheet so run.py can drive it: insert five
# metadata rows above the existing test-case table and re-point the
# selection/column-width to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: push the existing table down by 5 rows -----------------
$ws.Rows("1:5").Insert()
$ws.Rows("1:5").RowHeight = 12.8

# --- 2. Merge each metadata row's value area (B:F) first, so only the
#        anchor cell (column B) ends up materialised in sheetData ----------
$ws.Range("B1:F1").Merge()
$ws.Range("B2:F2").Merge()
$ws.Range("B3:F3").Merge()
$ws.Range("B4:F4").Merge()
$ws.Range("B5:F5").Merge()

# --- 3. Metadata labels in column A ----------------------------------------
$ws.Range("A1").Value = '$INC_FILE'
$ws.Range("A2").Value = '$FUNC_FILE'
$ws.Range("A3").Value = '$INC_PATH_GCC'
$ws.Range("A4").Value = '$FLAGS_GCC'
$ws.Range("A5").Value = '$FUNC_NAME'

# --- 4. Metadata values in column B (merge anchors) ------------------------
$ws.Range("B1").Value = "example/add.hpp"
$ws.Range("B2").Value = "example/add.cpp"
# B3 ($INC_PATH_GCC) intentionally left blank
$ws.Range("B4").Value = "-pthread"
$ws.Range("B5").Value = "add"

# --- 5. Alignment: left/center for text rows, center/center for the rest --
$ws.Range("B1").HorizontalAlignment = -4131
$ws.Range("B1").VerticalAlignment = -4108
$ws.Range("B2").HorizontalAlignment = -4131
$ws.Range("B2").VerticalAlignment = -4108
$ws.Range("B5").HorizontalAlignment = -4131
$ws.Range("B5").VerticalAlignment = -4108

$ws.Range("B3").HorizontalAlignment = -4108
$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("B4").HorizontalAlignment = -4108
$ws.Range("B4").VerticalAlignment = -4108

# --- 6. Widen column A slightly now that the labels are longer ------------
$ws.Columns("A").ColumnWidth = 16.53

# --- 7. Restore the expected selection -------------------------------------
$ws.Range("B5").Select() | Out-Null
